$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) sometimes contains values that look like plain
# numbers (e.g. "1.011"). Excel would normally auto-convert such text into
# a real number when assigned directly, which truncates meaningful trailing
# zeros (e.g. "7.170" -> 7.17) and changes the cell type. To preserve the
# exact original text representation, force those specific cells to Text
# format before writing the new value.

$ws.Range("D2").Value = "29.223.58"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.843.18"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.71"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6185"
$ws.Range("E6").Value = "  -2.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.010"
$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07443"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2948"
$ws.Range("E9").Value = "  +0.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.97"
$ws.Range("E10").Value = "  -0.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07717"
$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("D12").Value = "1.841.70"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.994"
$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6725"
$ws.Range("E14").Value = "  +0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.88"
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009101"
$ws.Range("E16").Value = "  -2.65%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.877"
$ws.Range("E17").Value = "  -2.89%  "

$ws.Range("D18").Value = "29.221.37"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "2.087.02"
$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.53"
$ws.Range("E20").Value = "  +5.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.61"
$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.011"
$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.170"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.015"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.63"
$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1424"
$ws.Range("E26").Value = "  +1.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.512"
$ws.Range("E27").Value = "  -0.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.88"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.503"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.143"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05581"
$ws.Range("E31").Value = "  +0.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.116"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.219"
$ws.Range("E33").Value = "  +0.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.851"
$ws.Range("E34").Value = "  -0.60%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7426"
$ws.Range("E35").Value = "  -1.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.140"
$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.656"
$ws.Range("E37").Value = "  +1.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.832"
$ws.Range("E38").Value = "  +2.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").Value = "  -0.72%  "

$ws.Range("D40").Value = "1.207.82"
$ws.Range("E40").Value = "  -2.50%  "

$ws.Range("E41").Value = "  -2.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9112"
$ws.Range("E42").Value = "  +1.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.59"
$ws.Range("E44").Value = "  -0.41%  "

$ws.Range("D45").Value = "1.991.38"
$ws.Range("E45").Value = "  +0.26%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000124"
$ws.Range("E46").Value = "  -0.96%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.15"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5137"
$ws.Range("E48").Value = "  +0.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.172"
$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4043"
$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05843"
$ws.Range("E51").Value = "  +0.40%  "

